# Auto-generated
# Parsed 94 changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.476.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.653.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '195.76'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '580.96'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.645.12'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.622'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.24'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +9.27%  '
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.154'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000298'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +18.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.224.43'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.644.37'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '68.366.37'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '405.74'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.84'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +25.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.26'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.46'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.72'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.89'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.10'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.20'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +22.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.25'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.86'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '693.42'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +17.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.30'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '43.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.423'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +13.20%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0802'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +9.12%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +20.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.15'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +13.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.215.13'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +19.29%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.98'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +31.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0424'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.30%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.92'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +9.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.13'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.18'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.77%  '
